$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.707.82"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").Value = "3.784.00"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "595.89"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("D6").Value = "167.32"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").Value = "3.782.57"
$ws.Range("E7").Value = "  +0.59%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -0.31%  "

$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("E11").Value = "  -1.95%  "

$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").Value = "0.0000253"
$ws.Range("E13").Value = "  -2.48%  "

$ws.Range("E14").Value = "  -0.11%  "

$ws.Range("D15").Value = "4.418.12"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").Value = "3.784.17"
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "18.44"
$ws.Range("E17").Value = "  +2.64%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "67.701.65"
$ws.Range("E18").Value = "  -1.17%  "

$ws.Range("E19").Value = "  +0.45%  "

$ws.Range("E20").Value = "  -1.11%  "

$ws.Range("D21").Value = "10.02"
$ws.Range("E21").Value = "  -6.40%  "

$ws.Range("D22").Value = "458.00"
$ws.Range("E22").Value = "  -1.55%  "

$ws.Range("E23").Value = "  -0.16%  "

$ws.Range("E24").Value = "  +4.58%  "

$ws.Range("E25").Value = "  -0.91%  "

$ws.Range("D26").Value = "11.96"
$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("E27").Value = "  -2.16%  "

$ws.Range("D28").Value = "10.05"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D31").Value = "2.25"
$ws.Range("E31").Value = "  +3.63%  "

$ws.Range("D32").Value = "7.22"
$ws.Range("E32").Value = "  -1.57%  "

$ws.Range("D34").Value = "9.11"
$ws.Range("E34").Value = "  -0.75%  "

$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("E36").Value = "  -0.24%  "

$ws.Range("D37").Value = "3.37"
$ws.Range("E37").Value = "  -0.63%  "

$ws.Range("E38").Value = "  +0.35%  "

$ws.Range("E39").Value = "  -0.53%  "

$ws.Range("E40").Value = "  -0.53%  "

$ws.Range("D41").Value = "1.00"

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("E43").Value = "  +4.36%  "

$ws.Range("D44").Value = "48.18"
$ws.Range("E44").Value = "  +3.14%  "

$ws.Range("E45").Value = "  -1.12%  "

$ws.Range("D46").Value = "149.46"
$ws.Range("E46").Value = "  +2.63%  "

$ws.Range("E47").Value = "  -1.89%  "

$ws.Range("D48").Value = "393.93"
$ws.Range("E48").Value = "  +1.08%  "

$ws.Range("D49").Value = "1.82"
$ws.Range("E49").Value = "  -4.71%  "

$ws.Range("D50").Value = "26.30"
$ws.Range("E50").Value = "  +1.20%  "

$ws.Range("D51").Value = "2.725.40"
$ws.Range("E51").Value = "  -1.27%  "
